$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typo "withing" -> "within" in the task description cell (C4)
$ws.Range("C4").Value = "Identify a subset use-case within the TOP use-case diagram"

# Update the active cell selection to C4 (was E26)
$ws.Range("C4").Select()
